# Update "想去人数" (want-to-go count) values in F column on sheets
# "展览" and "全部类型" for specific rows.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F10").Value = 429
    $ws.Range("F12").Value = 174
    $ws.Range("F26").Value = 4241
    $ws.Range("F32").Value = 613
    $ws.Range("F36").Value = 154
}
